$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row 19: Day 16, Date 17/3/2024, Time 1, Description
$ws.Range("A19").Value = 16
$ws.Range("B19").Value = "17/3/2024"
$ws.Range("C19").Value = 1
$ws.Range("D19").Value = "Added UpdateProtectedUserDetails functionality"

# Add new row 20: (no Day value), Date 19/3/2024, Time 3.16, Description
$ws.Range("B20").Value = "19/3/2024"
$ws.Range("C20").Value = 3.16
$ws.Range("D20").Value = "Added Product Category and Subcategory"

# Copy style from row 18 cells to the newly added row cells so formatting matches
$ws.Range("A18:D18").Copy()
$ws.Range("A19:D19").PasteSpecial(-4122)
$ws.Range("B18:D18").Copy()
$ws.Range("B20:D20").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Update the selected cell on the sheet (matches author's last-saved selection)
$ws.Range("D10").Select()
